$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B4").Value = 35238095.23809537
$ws.Range("C4").Value = 236628.9114398787
$ws.Range("E4").Value = 35238095.23809639

$ws.Range("C5").Value = 7887849.739434328

$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 21030415.97823774
$ws.Range("G6").Value = 27515499.01900803

$ws.Range("G7").Value = 524945.5733348924

# New row 13: Electrification + Bio-based feedstock
$ws.Range("A13").Value = "Electrification + Bio-based feedstock"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 28012865.72061013
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# New row 14: Conventional + Bio-based feedstock
$ws.Range("A14").Value = "Conventional + Bio-based feedstock"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 375614.1569601232
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

# New row 15: Conventional + Bio-based feedstock with CC
$ws.Range("A15").Value = "Conventional + Bio-based feedstock with CC"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 12520820.09792164
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0

# Copy style of existing label cell A12 to new label cells
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13:A15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wb.Save()
